$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = "Mark Vientos"
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 8
$ws.Range("H2").Value = 7
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 0.571
$ws.Range("K2").Value = 2
$ws.Range("A3").Value = "Luken Baker"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 20
$ws.Range("I3").Value = 0.333
$ws.Range("J3").Value = 0.05
$ws.Range("K3").Value = 0.333
$ws.Range("A4").Value = "Curtis Mead"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 18
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0.111
$ws.Range("K4").Value = 0
$ws.Range("A5").Value = "Hunter Feduccia"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 12
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 0.083
$ws.Range("K5").Value = 1
$ws.Range("A6").Value = "Andrew Benintendi"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 10
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0.1
$ws.Range("K6").Value = 0
$ws.Range("A7").Value = "Matt Thaiss"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3
$ws.Range("H7").Value = 10
$ws.Range("I7").Value = 0.667
$ws.Range("J7").Value = 0.2
$ws.Range("K7").Value = 1
$ws.Range("A8").Value = "Colby Thomas"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 11
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 0.091
$ws.Range("K8").Value = 2
$ws.Range("A9").Value = "Gio Urshela"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 12
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 0.083
$ws.Range("K9").Value = 1
$ws.Range("A10").Value = "Bryan Reynolds"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 17
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 0.059
$ws.Range("K10").Value = 1

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = "Tyrone Taylor"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 11
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 0.182
$ws.Range("K2").Value = 1
$ws.Range("A3").Value = "Griffin Conine"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 16
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 0.188
$ws.Range("K3").Value = 1
$ws.Range("A4").Value = "Marcelo Mayer"
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 13
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 0.308
$ws.Range("K4").Value = 1.5
$ws.Range("A5").Value = "Michael Conforto"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 19
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0.105
$ws.Range("K5").Value = 0
$ws.Range("A6").Value = "Alex Verdugo"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 13
$ws.Range("I6").Value = 0.5
$ws.Range("J6").Value = 0.077
$ws.Range("K6").Value = 0.5
$ws.Range("A7").Value = "Leo Jiménez"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 10
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 0.1
$ws.Range("K7").Value = 1
$ws.Range("A8").Value = "Jonathan Davis"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 9
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 0.222
$ws.Range("K8").Value = 1.5
$ws.Range("A9").Value = "Brandon Nimmo"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 10
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 0.2
$ws.Range("K9").Value = 1
$ws.Range("A10").Value = "Graham Pauley"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2
$ws.Range("H10").Value = 17
$ws.Range("I10").Value = 0.5
$ws.Range("J10").Value = 0.059
$ws.Range("K10").Value = 1

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = "Clarke Schmidt"
$ws.Range("B2").Value = 0.7
$ws.Range("C2").Value = 11
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 104
$ws.Range("J2").Value = 62
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = "104-62"
$ws.Range("A3").Value = "Jay Jackson"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = "8-2"
$ws.Range("A4").Value = "Luis Medina"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = "6-4"

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = "Mike Clevinger"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 98
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = "98-60"
$ws.Range("A3").Value = "J.P. France"
$ws.Range("B3").Value = 0.3
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 13
$ws.Range("J3").Value = 8
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = "13-8"
$ws.Range("A4").Value = "Ronny Henriquez"
$ws.Range("B4").Value = 0.3
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = "6-4"
